# Inclusion of beta distributions for compliance
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: tighten the (1-x) multiplier used to derive the NZ decline ratio ---
$ws.Range("F8").Formula = "=(1-0.85)*E8"

# --- Insert a fresh row at 11 so the date-block (old row 12 onward) shifts down by one ---
$ws.Rows("11:11").Insert()

# --- Row 9: new "Emissions reduction" input row (Aus/Wuhan style compliance ratio) ---
# G9 used to carry the "NZ distancing time %" label; it moves down to G10/G11 below.
$ws.Range("G9").ClearContents()
$ws.Range("B9").Value = "Emissions reduction"
$ws.Range("B9").Font.Bold = $true
$ws.Range("C9").Value = 0.414
$ws.Range("D9").Value = 0.283
$ws.Range("E9").Formula = "=D9/C9"
$ws.Range("E9").NumberFormat = "0%"
$ws.Range("F9").Formula = "=(1-0.85)*E9"
$ws.Range("F9").NumberFormat = "0.0%"

# --- Row 10: second "0.85 + F#" helper, carries over the label that used to sit on row 9 ---
$ws.Range("F10").Formula = "=0.85+F8"
$ws.Range("F10").NumberFormat = "0.0%"
$ws.Range("G10").Value = "NZ distancing time %"

# --- Row 11 (brand new, inserted above): matching helper for the new emissions-reduction row ---
$ws.Range("F11").Formula = "=0.85+F9"
$ws.Range("F11").NumberFormat = "0.0%"
$ws.Range("G11").Value = "NZ distancing time %"

# --- Row 12 (brand new, inserted above): empty but keeps the percent formatting of the column ---
$ws.Range("F12").NumberFormat = "0.0%"

# --- Row 18 (old row 17, shifted): new beta-distribution compliance ratio check ---
$ws.Range("G18").Formula = "=28.3/41.4 * 0.15"

# --- Restore the user's active selection ---
$ws.Range("G19").Select()
